$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, pushing existing rows 5-15 down to 6-16,
# copying formatting from the row above (matches the observed diff where the
# whole Primera/Segunda table gained one new weekly record at the top).
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new weekly price record.
$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(5, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(5, 4).Value = 45014
$ws.Cells.Item(5, 5).Value = 15
$ws.Cells.Item(5, 6).Value = "Fruta"
$ws.Cells.Item(5, 7).Value = 100107
$ws.Cells.Item(5, 8).Value = "Otros"
$ws.Cells.Item(5, 9).Value = 100107011
$ws.Cells.Item(5, 10).Value = "Tuna"
$ws.Cells.Item(5, 11).Value = "Sin especificar"
$ws.Cells.Item(5, 12).Value = "Segunda"
$ws.Cells.Item(5, 13).Value = 200
$ws.Cells.Item(5, 14).Value = 24000
$ws.Cells.Item(5, 15).Value = 25000
$ws.Cells.Item(5, 16).Value = 24500
$ws.Cells.Item(5, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(5, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(5, 19).Value = 1225
$ws.Cells.Item(5, 20).Value = 20
